$wb = $excel.ActiveWorkbook

# --- Update phone number value in Sheet1!C8 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C8").Value = "12345"

# --- Update selection/active cell on each sheet ---
$ws1.Range("C8").Select()

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B1").Select()

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A1").Select()

# Re-activate Sheet1 so it remains the tab that's selected/active
$ws1.Activate()
